$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.046.80"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "3.790.90"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'600.81"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("D6").Value = "'165.29"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("E9").Value = "  -0.97%  "

# Row 10
$ws.Range("D10").Value = "'0.451"
$ws.Range("E10").Value = "  +0.40%  "

# Row 11
$ws.Range("E11").Value = "  +2.63%  "

# Row 12
$ws.Range("E12").Value = "  -1.51%  "

# Row 13
$ws.Range("D13").Value = "'35.74"
$ws.Range("E13").Value = "  -0.96%  "

# Row 14
$ws.Range("D14").Value = "4.426.99"
$ws.Range("E14").Value = "  -0.31%  "

# Row 15
$ws.Range("D15").Value = "3.806.12"
$ws.Range("E15").Value = "  +0.77%  "

# Row 16
$ws.Range("D16").Value = "68.055.21"
$ws.Range("E16").Value = "  +0.39%  "

# Row 17
$ws.Range("E17").Value = "  -1.22%  "

# Row 18
$ws.Range("E18").Value = "  +1.93%  "

# Row 19
$ws.Range("E19").Value = "  -0.34%  "

# Row 20
$ws.Range("D20").Value = "'461.15"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
$ws.Range("D21").Value = "'9.71"
$ws.Range("E21").Value = "  -1.72%  "

# Row 22
$ws.Range("E22").Value = "  -0.51%  "

# Row 23
$ws.Range("D23").Value = "'0.0000149"
$ws.Range("E23").Value = "  -2.69%  "

# Row 24
$ws.Range("D24").Value = "'82.81"
$ws.Range("E24").Value = "  -0.84%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("E26").Value = "  +0.28%  "

# Row 27
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("D28").Value = "'9.99"
$ws.Range("E28").Value = "  -0.11%  "

# Row 29
$ws.Range("D29").Value = "3.940.79"
$ws.Range("E29").Value = "  -0.20%  "

# Row 30
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  +2.05%  "

# Row 31
$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  -5.42%  "

# Row 32
$ws.Range("E32").Value = "  -1.88%  "

# Row 33
$ws.Range("E33").Value = "  -1.30%  "

# Row 34
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").Value = "'8.99"
$ws.Range("E35").Value = "  -0.89%  "

# Row 36
$ws.Range("D36").Value = "'0.0999"
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("E37").Value = "  +0.45%  "

# Row 38
$ws.Range("E38").Value = "  -3.14%  "

# Row 39
$ws.Range("D39").Value = "'5.77"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("D40").Value = "'0.989"
$ws.Range("E40").Value = "  -0.63%  "

# Row 42
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
$ws.Range("D44").Value = "'47.37"
$ws.Range("E44").Value = "  -1.48%  "

# Row 45
$ws.Range("E45").Value = "  -2.34%  "

# Row 46
$ws.Range("D46").Value = "'151.44"
$ws.Range("E46").Value = "  +0.50%  "

# Row 47
$ws.Range("D47").Value = "'8.35"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48
$ws.Range("E48").Value = "  +2.76%  "

# Row 49
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'393.22"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.35"
$ws.Range("E50").Value = "  +6.50%  "

# Row 51
$ws.Range("D51").Value = "'26.77"
$ws.Range("E51").Value = "  +1.55%  "
